$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "'26.369.52"
$ws.Cells.Item(2, 5).Value = "'  -0.44%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "'1.830.82"
$ws.Cells.Item(3, 5).Value = "'  -0.55%  "

# Row 4
$ws.Cells.Item(4, 4).Value = "'1.001"
$ws.Cells.Item(4, 5).Value = "'  +0.06%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "'249.34"
$ws.Cells.Item(5, 5).Value = "'  -4.29%  "

# Row 6
$ws.Cells.Item(6, 5).Value = "'  +0.08%  "

# Row 7
$ws.Cells.Item(7, 4).Value = "'0.5246"
$ws.Cells.Item(7, 5).Value = "'  -0.14%  "

# Row 8
$ws.Cells.Item(8, 4).Value = "'0.2751"
$ws.Cells.Item(8, 5).Value = "'  -14.16%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "'0.06819"
$ws.Cells.Item(9, 5).Value = "'  +0.42%  "

# Row 10
$ws.Cells.Item(10, 4).Value = "'1.833.90"
$ws.Cells.Item(10, 5).Value = "'  -0.38%  "

# Row 11
$ws.Cells.Item(11, 4).Value = "'16.45"
$ws.Cells.Item(11, 5).Value = "'  -12.48%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "'0.07112"
$ws.Cells.Item(12, 5).Value = "'  -8.20%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "'0.6804"
$ws.Cells.Item(13, 5).Value = "'  -13.43%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "'85.66"
$ws.Cells.Item(14, 5).Value = "'  -2.35%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "'4.839"
$ws.Cells.Item(15, 5).Value = "'  -3.52%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "'1.002"
$ws.Cells.Item(16, 5).Value = "'  +0.09%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "'13.16"
$ws.Cells.Item(18, 5).Value = "'  -5.01%  "

# Row 19
$ws.Cells.Item(19, 2).Value = "'ShibaInu"
$ws.Cells.Item(19, 3).Value = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Cells.Item(19, 4).Value = "'0.000007309"
$ws.Cells.Item(19, 5).Value = "'  -7.98%  "

# Row 20
$ws.Cells.Item(20, 2).Value = "'WrappedBTC"
$ws.Cells.Item(20, 3).Value = "'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Cells.Item(20, 4).Value = "'26.391.52"
$ws.Cells.Item(20, 5).Value = "'  -0.45%  "

# Row 21
$ws.Cells.Item(21, 4).Value = "'2.073.99"
$ws.Cells.Item(21, 5).Value = "'  -0.08%  "

# Row 22
$ws.Cells.Item(22, 4).Value = "'4.479"
$ws.Cells.Item(22, 5).Value = "'  -3.32%  "

# Row 23
$ws.Cells.Item(23, 4).Value = "'5.791"
$ws.Cells.Item(23, 5).Value = "'  -3.24%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "'8.922"
$ws.Cells.Item(24, 5).Value = "'  -4.92%  "

# Row 25
$ws.Cells.Item(25, 4).Value = "'142.40"
$ws.Cells.Item(25, 5).Value = "'  +0.81%  "

# Row 26
$ws.Cells.Item(26, 4).Value = "'1.662"
$ws.Cells.Item(26, 5).Value = "'  -1.52%  "

# Row 27
$ws.Cells.Item(27, 4).Value = "'2.016"
$ws.Cells.Item(27, 5).Value = "'  -7.10%  "

# Row 28
$ws.Cells.Item(28, 4).Value = "'16.51"
$ws.Cells.Item(28, 5).Value = "'  -2.46%  "

# Row 29
$ws.Cells.Item(29, 4).Value = "'108.48"
$ws.Cells.Item(29, 5).Value = "'  -2.90%  "

# Row 30
$ws.Cells.Item(30, 4).Value = "'4.070"
$ws.Cells.Item(30, 5).Value = "'  -2.01%  "

# Row 31
$ws.Cells.Item(31, 4).Value = "'0.08700"
$ws.Cells.Item(31, 5).Value = "'  +0.15%  "

# Row 32
$ws.Cells.Item(32, 4).Value = "'3.833"
$ws.Cells.Item(32, 5).Value = "'  -5.87%  "

# Row 33
$ws.Cells.Item(33, 4).Value = "'0.04669"
$ws.Cells.Item(33, 5).Value = "'  -4.13%  "

# Row 34
$ws.Cells.Item(34, 4).Value = "'2.871"
$ws.Cells.Item(34, 5).Value = "'  +0.61%  "

# Row 35
$ws.Cells.Item(35, 4).Value = "'1.099"
$ws.Cells.Item(35, 5).Value = "'  -3.27%  "

# Row 36
$ws.Cells.Item(36, 4).Value = "'0.6947"
$ws.Cells.Item(36, 5).Value = "'  -4.73%  "

# Row 37
$ws.Cells.Item(37, 4).Value = "'3.043"
$ws.Cells.Item(37, 5).Value = "'  -1.52%  "

# Row 38
$ws.Cells.Item(38, 4).Value = "'2.167"
$ws.Cells.Item(38, 5).Value = "'  -3.71%  "

# Row 39
$ws.Cells.Item(39, 4).Value = "'0.01631"
$ws.Cells.Item(39, 5).Value = "'  -7.07%  "

# Row 40
$ws.Cells.Item(40, 4).Value = "'0.4437"
$ws.Cells.Item(40, 5).Value = "'  -7.06%  "

# Row 41
$ws.Cells.Item(41, 4).Value = "'0.8507"
$ws.Cells.Item(41, 5).Value = "'  -4.62%  "

# Row 42
$ws.Cells.Item(42, 4).Value = "'104.65"
$ws.Cells.Item(42, 5).Value = "'  -4.51%  "

# Row 43
$ws.Cells.Item(43, 5).Value = "'  +0.01%  "

# Row 44
$ws.Cells.Item(44, 4).Value = "'5.705"
$ws.Cells.Item(44, 5).Value = "'  -3.92%  "

# Row 45
$ws.Cells.Item(45, 4).Value = "'6.986"
$ws.Cells.Item(45, 5).Value = "'  -9.19%  "

# Row 46
$ws.Cells.Item(46, 4).Value = "'8.615"
$ws.Cells.Item(46, 5).Value = "'  -3.89%  "

# Row 47
$ws.Cells.Item(47, 2).Value = "'Elrond"
$ws.Cells.Item(47, 3).Value = "'https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Cells.Item(47, 4).Value = "'33.28"
$ws.Cells.Item(47, 5).Value = "'  -4.56%  "

# Row 48
$ws.Cells.Item(48, 2).Value = "'Cronos"
$ws.Cells.Item(48, 3).Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(48, 4).Value = "'0.05552"
$ws.Cells.Item(48, 5).Value = "'  -5.08%  "

# Row 49
$ws.Cells.Item(49, 2).Value = "'Aave"
$ws.Cells.Item(49, 3).Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(49, 4).Value = "'58.28"
$ws.Cells.Item(49, 5).Value = "'  -2.23%  "

# Row 50
$ws.Cells.Item(50, 4).Value = "'0.1169"
$ws.Cells.Item(50, 5).Value = "'  -5.18%  "

# Row 51
$ws.Cells.Item(51, 4).Value = "'0.8555"
$ws.Cells.Item(51, 5).Value = "'  -4.22%  "
